# Apply the "assess by word" update: add a new word (ਗੁਬਾਰੀ) entry to the
# "Words" sheet and its three matching verse rows to the "Progress" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Words": append row 3 with the new word, mirroring row 2's layout.
# ---------------------------------------------------------------------
$wsWords = $wb.Worksheets.Item("Words")

$wsWords.Cells.Item(3, 1).Value = "ਗੁਬਾਰੀ"          # A3 word
$wsWords.Cells.Item(3, 2).Value = "ਗੁਬਾਰੀ"          # B3 word_key_norm
$wsWords.Cells.Item(3, 3).Value = $true              # C3 listed_by_user

$wsWords.Cells.Item(3, 4).Value = 45914.46202929398  # D3 listed_at
$wsWords.Cells.Item(3, 4).NumberFormat = $wsWords.Cells.Item(2, 4).NumberFormat

$wsWords.Cells.Item(3, 5).Value = $true              # E3 selected_for_analysis

$wsWords.Cells.Item(3, 6).Value = 45914.46202929398  # F3 selected_at
$wsWords.Cells.Item(3, 6).NumberFormat = $wsWords.Cells.Item(2, 6).NumberFormat

$wsWords.Cells.Item(3, 7).Value = $true              # G3 analysis_started

$wsWords.Cells.Item(3, 8).Value = 45914.46206456018  # H3 analysis_started_at
$wsWords.Cells.Item(3, 8).NumberFormat = $wsWords.Cells.Item(2, 8).NumberFormat

$wsWords.Cells.Item(3, 9).Value = $false             # I3 analysis_completed
$wsWords.Cells.Item(3, 10).Value = ""                # J3 analysis_completed_at
$wsWords.Cells.Item(3, 11).Value = 0                 # K3 sequence_index
$wsWords.Cells.Item(3, 12).Value = ""                # L3 notes

# ---------------------------------------------------------------------
# Sheet "Progress": append rows 5-7, one per verse containing the word.
# ---------------------------------------------------------------------
$wsProgress = $wb.Worksheets.Item("Progress")

# Row 5 - page 507
$wsProgress.Cells.Item(5, 1).Value = "ਗੁਬਾਰੀ"
$wsProgress.Cells.Item(5, 2).Value = "ਗੁਬਾਰੀ"
$wsProgress.Cells.Item(5, 3).Value = ""
$wsProgress.Cells.Item(5, 4).Value = "ਮਨਮੁਖ ਦੁਬਿਧਾ ਦੁਰਮਤਿ ਬਿਆਪੇ ਜਿਨ ਅੰਤਰਿ ਮੋਹ ਗੁਬਾਰੀ ॥"
$wsProgress.Cells.Item(5, 5).Value = 507
$wsProgress.Cells.Item(5, 6).Value = $true
$wsProgress.Cells.Item(5, 7).Value = 45914.46202929398
$wsProgress.Cells.Item(5, 7).NumberFormat = $wsProgress.Cells.Item(2, 7).NumberFormat
$wsProgress.Cells.Item(5, 8).Value = "not started"
$wsProgress.Cells.Item(5, 9).Value = ""
$wsProgress.Cells.Item(5, 10).Value = ""
$wsProgress.Cells.Item(5, 11).Value = ""

# Row 6 - page 788
$wsProgress.Cells.Item(6, 1).Value = "ਗੁਬਾਰੀ"
$wsProgress.Cells.Item(6, 2).Value = "ਗੁਬਾਰੀ"
$wsProgress.Cells.Item(6, 3).Value = ""
$wsProgress.Cells.Item(6, 4).Value = "ਜਿਨ੍ਹ੍ਹਿ ਕੀਏ ਤਿਸਹਿ ਨ ਜਾਣਨੀ ਮਨਮੁਖਿ ਗੁਬਾਰੀ ॥"
$wsProgress.Cells.Item(6, 5).Value = 788
$wsProgress.Cells.Item(6, 6).Value = $true
$wsProgress.Cells.Item(6, 7).Value = 45914.46202929398
$wsProgress.Cells.Item(6, 7).NumberFormat = $wsProgress.Cells.Item(2, 7).NumberFormat
$wsProgress.Cells.Item(6, 8).Value = "not started"
$wsProgress.Cells.Item(6, 9).Value = ""
$wsProgress.Cells.Item(6, 10).Value = ""
$wsProgress.Cells.Item(6, 11).Value = ""

# Row 7 - page 1243
$wsProgress.Cells.Item(7, 1).Value = "ਗੁਬਾਰੀ"
$wsProgress.Cells.Item(7, 2).Value = "ਗੁਬਾਰੀ"
$wsProgress.Cells.Item(7, 3).Value = ""
$wsProgress.Cells.Item(7, 4).Value = "ਬਾਹਰਿ ਭਸਮ ਲੇਪਨ ਕਰੇ ਅੰਤਰਿ ਗੁਬਾਰੀ ॥"
$wsProgress.Cells.Item(7, 5).Value = 1243
$wsProgress.Cells.Item(7, 6).Value = $true
$wsProgress.Cells.Item(7, 7).Value = 45914.46202929398
$wsProgress.Cells.Item(7, 7).NumberFormat = $wsProgress.Cells.Item(2, 7).NumberFormat
$wsProgress.Cells.Item(7, 8).Value = "not started"
$wsProgress.Cells.Item(7, 9).Value = ""
$wsProgress.Cells.Item(7, 10).Value = ""
$wsProgress.Cells.Item(7, 11).Value = ""
